# Rows 15, 17 and 18 each get the record that used to live in the "next"
# row in the 15 -> 17 -> 18 -> 15 cycle:
#   new row15 = old row17
#   new row17 = old row18
#   new row18 = old row15
# Row 16 (and every other row) is untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 15  (<- old Row 17: Grönpyrola / Pyrola chlorantha) ----
$ws.Range("A15").Value = 111837675
$ws.Range("B15").Value = 103288
$ws.Range("D15").Value = "LC"
$ws.Range("E15").Value = 221144
$ws.Range("F15").Value = "Grönpyrola"
$ws.Range("G15").Value = "Pyrola chlorantha"
$ws.Range("H15").Value = "Sw."
$ws.Range("I15").Value = "'10"
$ws.Range("I15").Style = "Normal"
$ws.Range("J15").Value = "plantor/tuvor"
$ws.Range("L15").Style = "Normal"
$ws.Range("P15").Value = "Brotorp, Långsjön, Sm"
$ws.Range("Q15").Value = 575781.9606960951
$ws.Range("R15").Value = 6404546.96767282
$ws.Range("AC15").ClearContents()

# ---- Row 17  (<- old Row 18: Zontaggsvamp / Hydnellum concrescens) ----
$ws.Range("A17").Value = 111837705
$ws.Range("B17").Value = 90662
$ws.Range("E17").Value = 4363
$ws.Range("F17").Value = "Zontaggsvamp"
$ws.Range("G17").Value = "Hydnellum concrescens"
$ws.Range("H17").Value = "(Pers.) Banker"
$ws.Range("J17").Value = "fruktkroppar"
$ws.Range("L17").ClearContents()
$ws.Range("Q17").Value = 575795.3141537429
$ws.Range("R17").Value = 6404518.948622406

# ---- Row 18  (<- old Row 15: Koralltaggsvamp / Hericium coralloides) ----
$ws.Range("A18").Value = 111837758
$ws.Range("B18").Value = 90187
$ws.Range("D18").Value = "NT"
$ws.Range("E18").Value = 2014
$ws.Range("F18").Value = "Koralltaggsvamp"
$ws.Range("G18").Value = "Hericium coralloides"
$ws.Range("H18").Value = "(Scop.:Fr.) Pers."
$ws.Range("I18").Value = "'6"
$ws.Range("I18").Style = "Normal"
$ws.Range("P18").Value = "Brotorp, hyggeskant, Sm"
$ws.Range("Q18").Value = 575673.5681218
$ws.Range("R18").Value = 6404513.458820416
$ws.Range("AC18").Value = "På asplåga."
